$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.893646240234375
$ws.Range("B1").Value = 1.589063048362732
$ws.Range("C1").Value = 6.080820083618164
$ws.Range("D1").Value = 1.844580411911011
$ws.Range("E1").Value = 1.111282587051392
